# Apply data repull / push / mean-calculation updates to column F (dSF)
# on Sheet1. Only the "dSF" values for a handful of rows changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "F2"  = -6
    "F4"  = -3
    "F7"  = 0
    "F10" = -9
    "F12" = -5
    "F15" = -1
    "F22" = 0
    "F27" = -3
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
